$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.74313608499204076
$ws.Range("E1").Value = 0.56536580948380344
$ws.Range("V1").Value = 0.69041005924920229
$ws.Range("BJ1").Value = 0.98118328981227154
$ws.Range("M2").Value = 0.80008323777336865
$ws.Range("N3").Value = 0.81237450520765198
$ws.Range("BB3").Value = 0.99806341687359035
$ws.Range("E4").Value = 0.99115169074890863
$ws.Range("V4").Value = 0.9158183680731049
$ws.Range("BK5").Value = 0.85965722511811404
$ws.Range("K7").Value = 0.86940810182621575
$ws.Range("AY7").Value = 0.92368428406678593
$ws.Range("L8").Value = 0.85123647095088784
$ws.Range("V8").Value = 0.85729062000626044
$ws.Range("AR8").Value = 0.80112484509019333
$ws.Range("L10").Value = 0.8330522142437401
$ws.Range("Y10").Value = 0.90703790630522096
$ws.Range("AA10").Value = 0.88025587314424081
$ws.Range("AU10").Value = 0.72730545343167075
$ws.Range("D11").Value = 0.80198541782281407
$ws.Range("AL12").Value = 0.74773082790098711
$ws.Range("I13").Value = 0.70177446764466722
$ws.Range("Q13").Value = 0.98786733439634822
$ws.Range("AJ13").Value = 0.79747480340629817
$ws.Range("AK13").Value = 0.72248321870761623
$ws.Range("AT13").Value = 0.95599023977272557
$ws.Range("BP13").Value = 0.64514934139993207
$ws.Range("BF14").Value = 0.95148364487046266
$ws.Range("N17").Value = 0.83986445102276286
$ws.Range("P17").Value = 0.69733574741330795
$ws.Range("V18").Value = 0.98181471047235092
$ws.Range("AC18").Value = 0.9154585480550983
$ws.Range("AA20").Value = 0.86843099346119756
$ws.Range("BI20").Value = 0.80414002261487405
$ws.Range("E21").Value = 0.76606715786392388
$ws.Range("W21").Value = 0.93539002729111342
$ws.Range("AX21").Value = 0.65543904172278811
$ws.Range("BD21").Value = 0.80959521433646953
$ws.Range("AD22").Value = 0.79917959871674471
$ws.Range("BJ22").Value = 0.81697071655984754
$ws.Range("Z24").Value = 0.60309862927422864
$ws.Range("AV24").Value = 0.92675833473331415
$ws.Range("S25").Value = 0.80364506935476532
$ws.Range("AM26").Value = 0.85117183929857987
$ws.Range("N27").Value = 0.94095543709559415
$ws.Range("AE27").Value = 0.99074084763046288
$ws.Range("AL27").Value = 0.68033256232679085
$ws.Range("BF27").Value = 0.76515647726864744
$ws.Range("AH28").Value = 0.5802046358409273
$ws.Range("P29").Value = 0.88918122495868401
$ws.Range("AB29").Value = 0.94181490958186387
$ws.Range("AS29").Value = 0.98492320502886499
$ws.Range("BJ29").Value = 0.89245150680293173
$ws.Range("N30").Value = 0.9051544707556759
$ws.Range("P30").Value = 0.73135118714557867
$ws.Range("AK30").Value = 0.97895565887007896
$ws.Range("U32").Value = 0.89109608177281285
$ws.Range("X32").Value = 0.83460211007016194
$ws.Range("AK32").Value = 0.83027068312663532
$ws.Range("AP32").Value = 0.99387310058384026
$ws.Range("BN32").Value = 0.91469168328039951
$ws.Range("AE33").Value = 0.71092525583303523
$ws.Range("AH33").Value = 0.59616578078752913
$ws.Range("BL34").Value = 0.95023912708241887
$ws.Range("R35").Value = 0.92620084763253874
$ws.Range("T35").Value = 0.9372112372288266
$ws.Range("AG35").Value = 0.7051880782500477
$ws.Range("AU35").Value = 0.75250965353585419
$ws.Range("K37").Value = 0.91811297683924975
$ws.Range("W38").Value = 0.98385375198927572
$ws.Range("AO38").Value = 0.83227482352318405
$ws.Range("A39").Value = 0.82797711191323697
$ws.Range("AI39").Value = 0.7970086882241979
$ws.Range("N40").Value = 0.90602254115415792
$ws.Range("X41").Value = 0.5600394915564415
$ws.Range("Y41").Value = 0.83512608846959646
$ws.Range("AK41").Value = 0.62384862187442702
$ws.Range("AM41").Value = 0.80388040833998242
$ws.Range("O42").Value = 0.90451469557783282
$ws.Range("AC42").Value = 0.97134081830867147
$ws.Range("BM42").Value = 0.73643443732097502
$ws.Range("AH43").Value = 0.72192677776569836
$ws.Range("AT43").Value = 0.8052203144487925
$ws.Range("AY43").Value = 0.99494607030714277
$ws.Range("AS44").Value = 0.99733242913753584
$ws.Range("AJ45").Value = 0.89526485103708298
$ws.Range("Q46").Value = 0.90808159593845972
$ws.Range("BB47").Value = 0.8163446256834801
$ws.Range("I48").Value = 0.59031572607854277
$ws.Range("L48").Value = 0.59956040981701508
$ws.Range("BB48").Value = 0.82776191897370222
$ws.Range("D49").Value = 0.78642651342466918
$ws.Range("O49").Value = 0.71364519307832308
$ws.Range("AJ49").Value = 0.68312369504757531
$ws.Range("BF50").Value = 0.71103382268252635
$ws.Range("C52").Value = 0.84410566978133095
$ws.Range("T52").Value = 0.9907503054867941
$ws.Range("AN52").Value = 0.6574810332830261
$ws.Range("BF52").Value = 0.81165848433272303
$ws.Range("W53").Value = 0.75280551669852191
$ws.Range("AN53").Value = 0.69892208045457149
$ws.Range("BP53").Value = 0.94799752397784731
$ws.Range("A54").Value = 0.56126316181239022
$ws.Range("X54").Value = 0.65598489296712892
$ws.Range("V55").Value = 0.82477677474606348
$ws.Range("AL55").Value = 0.68493276864483665
$ws.Range("AS55").Value = 0.77137685236483899
$ws.Range("BH55").Value = 0.70722011373640059
$ws.Range("AH56").Value = 0.63966277614223244
$ws.Range("AR56").Value = 0.91886078228089407
$ws.Range("AT56").Value = 0.84086225307732765
$ws.Range("S57").Value = 0.99509760704204031
$ws.Range("AB57").Value = 0.9653548853643551
$ws.Range("F58").Value = 0.79052418733727414
$ws.Range("BE58").Value = 0.85433330953556186
$ws.Range("H59").Value = 0.54388199685214333
$ws.Range("BH59").Value = 0.98449157295030942
$ws.Range("F60").Value = 0.89802538701520318
$ws.Range("AO60").Value = 0.8626856664703858
$ws.Range("BO60").Value = 0.92783211565454038
$ws.Range("AM61").Value = 0.84228400290325067
$ws.Range("Z63").Value = 0.92939651401379964
$ws.Range("AC63").Value = 0.74606689396163395
$ws.Range("AD63").Value = 0.98599920461972101
$ws.Range("AV63").Value = 0.99343749441433826
$ws.Range("AA64").Value = 0.83701457405391544
$ws.Range("BN64").Value = 0.87295659899744238
$ws.Range("AH65").Value = 0.89148779189832106
$ws.Range("AY66").Value = 0.91827866376990563
$ws.Range("B67").Value = 0.8716787027593691
$ws.Range("F67").Value = 0.99506994792473558
$ws.Range("BK67").Value = 0.86401003029878654
$ws.Range("BN67").Value = 0.75454889086510335
$ws.Range("P68").Value = 0.87701426959868589
$ws.Range("AB68").Value = 0.58622123079971833
$ws.Range("AF68").Value = 0.95414669999341717
